$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.999.54'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').Value = '1.641.63'
$ws.Range('E3').Value = '  +0.47%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  +0.41%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '215.80'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -0.13%  '
$ws.Range('E8').Value = '  +0.36%  '
$ws.Range('E9').Value = '  +0.65%  '
$ws.Range('E10').Value = '  -0.38%  '
$ws.Range('E11').Value = '  +0.37%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '4.27'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +0.50%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.633.13'
$ws.Range('E13').Value = '  -0.24%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.544'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -0.15%  '
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').Value = '0.0₃0762'
$ws.Range('E15').Value = '  +0.72%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '63.37'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +1.26%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '26.036.11'
$ws.Range('E17').Value = '  +0.61%  '
$ws.Range('B18').Value = 'Dai'
$ws.Range('C18').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '1.01'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +0.37%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '194.01'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '4.36'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -0.77%  '
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '9.92'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('B22').Value = 'Chainlink'
$ws.Range('C22').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '6.20'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -1.27%  '
$ws.Range('B23').Value = 'Stellar'
$ws.Range('C23').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '0.132'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +4.74%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '1.80'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -1.16%  '
$ws.Range('B25').Value = 'BinanceUSD'
$ws.Range('C25').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '1.01'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +0.36%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '142.95'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -0.36%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '6.87'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +0.42%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '15.52'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +0.57%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '1.25'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +0.53%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '0.0495'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -0.91%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '3.28'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -0.45%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '3.26'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +1.13%  '
$ws.Range('B33').Value = 'LidoDAOToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '1.54'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -1.22%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '2.46'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +1.28%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.904'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +0.28%  '
$ws.Range('B36').Value = 'Maker'
$ws.Range('C36').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D36').Value = '1.130.25'
$ws.Range('E36').Value = '  -0.78%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.539'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -1.02%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '2.46'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -0.47%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.0157'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +0.17%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '5.47'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +0.91%  '
$ws.Range('B41').Value = 'Quant'
$ws.Range('C41').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '99.02'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -0.35%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.797'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -0.28%  '
$ws.Range('B43').Value = 'RocketPoolETH'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D43').Value = '1.778.05'
$ws.Range('E43').Value = '  +0.58%  '
$ws.Range('B44').Value = 'BabyDogeCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D44').Value = '0.0₆0117'
$ws.Range('E44').Value = '  +4.67%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '56.56'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +0.45%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.0522'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -1.28%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '1.49'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +2.92%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '7.73'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +1.01%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.414'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -0.18%  '
$ws.Range('B50').Value = 'USDD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +0.29%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.0952'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -0.78%  '
